$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-23 (Num, Error) - refreshed raw data
$data = @(
    @("КПМ-0273", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0283", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0285", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0286", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0287", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0289", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0291", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0292", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0293", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0294", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0295", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0296", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0341", "['Не выбрано ни одной проблемы']"),
    @("КПМ-0363", "['Не выбрано ни одной проблемы']"),
    @("Лаб-2022-10", "['Не выбрано ни одной проблемы', 'Не выбрано ни одной проблемы']"),
    @("Лаб-2022-12", "['Не выбрано ни одной проблемы']"),
    @("Лаб-2022-14", "['Не выбрано ни одной проблемы', 'Не выбрано ни одной проблемы', 'Не выбрано ни одной проблемы']"),
    @("Лаб-2022-17", "['Не выбрано ни одной проблемы']"),
    @("Лаб-2022-18", "['Не выбрано ни одной проблемы']"),
    @("Лаб-2022-2", "['Не выбрано ни одной проблемы', 'Не выбрано ни одной проблемы']"),
    @("Лаб-2022-21", "['Не выбрано ни одной проблемы']"),
    @("Лаб-2022-34", "['Не выбрано ни одной проблемы']")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}

# The new data only spans down to row 23; delete the now-stale rows 24-33
$ws.Range("A24:B33").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
